$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph in "Section 1: Executive Summary" whose
#    entire text is just "as" and replace it (plus add two new
#    paragraphs after it) with the full executive-summary copy.
#    We build the exact target OOXML for the three paragraphs
#    (including paragraph-mark rPr / bookmark placement) and drop it
#    in via Range.InsertXML so formatting round-trips precisely.
# ------------------------------------------------------------------

$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "as") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'as' placeholder paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)
$targetRange = $target.Range

$xmlFragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="9" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>The project will focus on the creation of a comprehensive web-based application for supporting migrants in need, with special emphasis on job opportunities, basic services, law and culture. New immigrants may face numerous issues in the short term upon arrival, such as language barriers, social isolation and lack of job opportunities, all of which can pose challenges to their basic livelihoods.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>This web application aims to provide a platform for concerned organizations and individuals to seek support and advice, and to help them find or give access to key resources such as health education, social services and job opportunities. With a wide range of connections between immigrants and society, this platform will serve immigrants to enjoy their entitlements while contributing to this country.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>It will target different audiences, including new immigrants, refugees and relevant volunteer or charitable organizations. The main goal is to meet the needs of migrants by providing content and services tailored to their specific needs, ultimately making them easily accessible even to those with limited literacy or language skills.</w:t></w:r><w:bookmarkEnd w:id="9"/></w:p>
'@

$targetRange.InsertXML($xmlFragment)

# ------------------------------------------------------------------
# 2. Style-sheet tweaks: mark "heading 4" and "footer" as quick
#    (qFormat) styles, matching the other heading/built-in styles.
# ------------------------------------------------------------------

$heading4 = $d.Styles("heading 4")
$heading4.QuickStyle = $true

$footer = $d.Styles("footer")
$footer.QuickStyle = $true

Write-Output "edit complete"
